{"js": "// Update the 25 multiplication problems (5 columns x 5 populated rows) in\n// the single table of the document. Cells are addressed by their\n// (row, column) position in the table rather than by searching for the old\n// text, because some of the new values collide with other cells' old\n// values (e.g. \"804\u00d72=\" is both an old value in one cell and a new value\n// written into a different cell), which would make a naive global\n// find-and-replace order-dependent and unsafe.\nconst table = context.document.body.tables.getFirst();\ntable.load(\"rowCount\");\nawait context.sync();\n\n// Map of table-row-index -> [new value for each of the 5 columns].\n// Only every 5th row (0, 4, 9, 14, 19) actually holds a problem; the rows\n// in between are empty spacer rows and are left untouched.\nconst newValuesByRow = {\n  0: [\"466\u00d77=\", \"468\u00d76=\", \"475\u00d79=\", \"104\u00d76=\", \"882\u00d72=\"],\n  4: [\"804\u00d72=\", \"668\u00d79=\", \"166\u00d72=\", \"455\u00d78=\", \"415\u00d76=\"],\n  9: [\"135\u00d77=\", \"283\u00d74=\", \"143\u00d76=\", \"412\u00d76=\", \"925\u00d75=\"],\n  14: [\"999\u00d78=\", \"783\u00d76=\", \"219\u00d73=\", \"980\u00d74=\", \"695\u00d77=\"],\n  19: [\"734\u00d77=\", \"425\u00d76=\", \"911\u00d74=\", \"977\u00d74=\", \"790\u00d75=\"],\n};\n\nfor (const rowIndexStr of Object.keys(newValuesByRow)) {\n  const rowIndex = Number(rowIndexStr);\n  const rowValues = newValuesByRow[rowIndex];\n  for (let col = 0; col < rowValues.length; col++) {\n    const cell = table.getCell(rowIndex, col);\n    cell.value = rowValues[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the 25 multiplication problems (5 columns x 5 populated rows) in\n# the single table of the document. Cells are addressed by their\n# (row, column) position in the table rather than by searching for the old\n# text, because some of the new values collide with other cells' old\n# values (e.g. \"804x2=\" is both an old value in one cell and a new value\n# written into a different cell), which would make a naive global\n# find-and-replace order-dependent and unsafe.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Word table rows/columns are 1-based. Only every 5th row (1, 5, 10, 15, 20)\n# actually holds a problem; the rows in between are empty spacer rows and\n# are left untouched.\n$newValues = @{\n    1  = @(\"466\u00d77=\", \"468\u00d76=\", \"475\u00d79=\", \"104\u00d76=\", \"882\u00d72=\")\n    5  = @(\"804\u00d72=\", \"668\u00d79=\", \"166\u00d72=\", \"455\u00d78=\", \"415\u00d76=\")\n    10 = @(\"135\u00d77=\", \"283\u00d74=\", \"143\u00d76=\", \"412\u00d76=\", \"925\u00d75=\")\n    15 = @(\"999\u00d78=\", \"783\u00d76=\", \"219\u00d73=\", \"980\u00d74=\", \"695\u00d77=\")\n    20 = @(\"734\u00d77=\", \"425\u00d76=\", \"911\u00d74=\", \"977\u00d74=\", \"790\u00d75=\")\n}\n\nforeach ($row in $newValues.Keys) {\n    $rowValues = $newValues[$row]\n    for ($col = 1; $col -le $rowValues.Count; $col++) {\n        $cell = $t.Cell($row, $col)\n        $cell.Range.Text = $rowValues[$col - 1]\n    }\n}\n"}
